$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -1.478419188821007
$ws.Range("B2").Value = -4.707344516337067

$ws.Range("A3").Value = -0.4858314172122379
$ws.Range("B3").Value = 0.5802011976966502

$ws.Range("A4").Value = 0.863089517208964
$ws.Range("B4").Value = -2.957586405059521

$ws.Range("B5").Value = 0.5424457742420091

$ws.Range("A6").Value = -0.8161408853368732
$ws.Range("B6").Value = -1.972806165615633

$ws.Range("A7").Value = -0.07197838678681279
$ws.Range("B7").Value = -0.6742525705485135

$ws.Range("A8").Value = 0.7793861571261941
$ws.Range("B8").Value = 0.9198368732398564

$ws.Range("A9").Value = 0.3188066342661657
$ws.Range("B9").Value = -0.2078810524838768

$ws.Range("A10").Value = -0.1811713301464574
$ws.Range("B10").Value = -1.82992045462128
